# Fix the "nan"/float formatting issue in the FRIDAY (column F) schedule.
# The source data had been round-tripped through something that turned
# integer period numbers into floats ("2.0", "4.0", ...) and missing
# cells into the literal text "nan". This restores the FRIDAY column to
# match the integer-style labels used by the other weekday columns, and
# clears the cells that should have been empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that are purely numeric-looking ("2.0" -> "2", etc.) need to be
# forced to text so Excel keeps storing them as shared strings (matching
# the other weekday columns) instead of re-typing them as numbers.
$plainNumericCells = @("F2", "F4", "F5", "F6", "F7")

foreach ($addr in $plainNumericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("F2").Replace("2.0", "2")
$ws.Range("F4").Replace("6.0", "6")
$ws.Range("F5").Replace("8.0", "8")
$ws.Range("F6").Replace("10.0", "10")
$ws.Range("F7").Replace("12.0", "12")

# Restore the cell formatting so we don't leave a stray text format behind.
foreach ($addr in $plainNumericCells) {
    $ws.Range($addr).ClearFormats()
}

# F3 keeps its trailing course info, so Replace naturally preserves the
# string type (the cell never looks like a pure number).
$ws.Range("F3").Replace("4.0", "4")

# F8:F10 held the literal text "nan" - these should simply be empty.
$ws.Range("F8:F10").ClearContents()
